$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# --- Insert a new physical row at 89, shifting everything below down ---
$ws.Rows("89:89").Insert()

# Copy formatting from the row below (old row 89, now row 90) into the
# newly-blanked row 89 so styles match the surrounding table rows.
$ws.Range("A90:K90").Copy() | Out-Null
$ws.Range("A89:K89").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the calculated-column formula in the new row's "EARNED " column
$ws.Range("G89").Formula = '=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"",Table13[[#This Row],[EARNED]])'

# --- Resize the table so it covers the newly inserted row too ---
$lo = $ws.ListObjects.Item("Table13")
$lo.Resize($ws.Range("A8:K136")) | Out-Null

# Make sure the calculated-column formula on the table's new last row
# (shifted from 135 -> 136) is the full structured-reference form.
$ws.Range("G136").Formula = '=IF(ISBLANK(Table13[[#This Row],[EARNED]]),"",Table13[[#This Row],[EARNED]])'

# --- Row 87: SL(3-0-0) taken 8/3,4,7/2023 earns an extra 1.25 VL/SL credit ---
$ws.Range("C87").Value = 1.25

# --- Row 88: SL(3-0-0) leave particulars, 1.25 credit, 3 days absence (SL) ---
$ws.Range("B88").Value = "SL(3-0-0)"
$ws.Range("C88").Value = 1.25
$ws.Range("H88").Value = 3
$ws.Range("K88").Value = "8/3,4,7/2023"

# --- Row 89 (new): SL(2-0-0) taken 8/24,25/2023, 2 days absence (SL) ---
$ws.Range("B89").Value = "SL(2-0-0)"
$ws.Range("H89").Value = 2
$ws.Range("K89").Value = "8/24,25/2023"
